# TC13_Canine_Filter_Breed-Bulldog.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab Cypher query (cell B2 on the "startup" sheet) had an
# erroneous trailing `Cohort` column (coalesce(co.cohort_description, '')
# AS `Cohort`) that was causing query errors. Remove that trailing
# column from the RETURN clause so the query matches the other,
# already-correct tabs (SamplesTab / FilesTab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bulldog']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Author re-selected B2 (previously B4 was selected, with the view
# scrolled so row 4 was the top row) before saving.
$ws.Activate()
$ws.Range("B2").Select()
